$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# The original sheet has 31 data rows (rows 1-31). The updated report restructures
# the "LONG XUYEN" and "SOC TRANG" sections to each include three additional lines
# ("Tong cong", "Phu cap", "Luong cong tac"), growing the sheet to 37 rows.
# Insert 3 blank rows before row 12 (start of the LONG XUYEN block).
$ws.Rows.Item(12).Resize(3, 1).EntireRow.Insert()

# Insert 3 blank rows before row 23 (start of the SOC TRANG block, after the first insert).
$ws.Rows.Item(23).Resize(3, 1).EntireRow.Insert()

# Write every label/value pair for the final layout (rows 1-37).
$ws.Range("A1").Value2 = "Danh mục lương"
$ws.Range("B1").Value2 = 3
$ws.Range("A2").Value2 = "Tổng công tại CẦN THƠ"
$ws.Range("B2").Value2 = 28
$ws.Range("A3").Value2 = "Phụ cấp tại CẦN THƠ"
$ws.Range("B3").Value2 = 980000
$ws.Range("A4").Value2 = "Lương cơ bản tại CẦN THƠ"
$ws.Range("B4").Value2 = 6000000
$ws.Range("A5").Value2 = "Chiết khấu sale chính tại CẦN THƠ"
$ws.Range("B5").Value2 = 0
$ws.Range("A6").Value2 = "Chiết khấu sale phụ tại CẦN THƠ"
$ws.Range("B6").Value2 = 0
$ws.Range("A7").Value2 = "Đơn 1 bác sĩ tại CẦN THƠ"
$ws.Range("B7").Value2 = 0
$ws.Range("A8").Value2 = "Đơn 2 bác sĩ tại CẦN THƠ"
$ws.Range("B8").Value2 = 0
$ws.Range("A9").Value2 = "Công phụ phẫu 1 tại CẦN THƠ"
$ws.Range("B9").Value2 = 0
$ws.Range("A10").Value2 = "Công phụ phẫu 2 tại CẦN THƠ"
$ws.Range("B10").Value2 = 0
$ws.Range("A11").Value2 = "Ứng lương tại CẦN THƠ"
$ws.Range("B11").Value2 = 0
$ws.Range("A12").Value2 = "Tổng công tại LONG XUYÊN"
$ws.Range("B12").Value2 = 0
$ws.Range("A13").Value2 = "Phụ cấp tại LONG XUYÊN"
$ws.Range("B13").Value2 = 0
$ws.Range("A14").Value2 = "Lương công tác tại LONG XUYÊN"
$ws.Range("B14").Value2 = 0
$ws.Range("A15").Value2 = "Lương cơ bản tại LONG XUYÊN"
$ws.Range("B15").ClearContents()
$ws.Range("A16").Value2 = "Chiết khấu sale chính tại LONG XUYÊN"
$ws.Range("B16").Value2 = 0
$ws.Range("A17").Value2 = "Chiết khấu sale phụ tại LONG XUYÊN"
$ws.Range("B17").Value2 = 0
$ws.Range("A18").Value2 = "Đơn 1 bác sĩ tại LONG XUYÊN"
$ws.Range("B18").Value2 = 0
$ws.Range("A19").Value2 = "Đơn 2 bác sĩ tại LONG XUYÊN"
$ws.Range("B19").Value2 = 0
$ws.Range("A20").Value2 = "Công phụ phẫu 1 tại LONG XUYÊN"
$ws.Range("B20").Value2 = 0
$ws.Range("A21").Value2 = "Công phụ phẫu 2 tại LONG XUYÊN"
$ws.Range("B21").Value2 = 0
$ws.Range("A22").Value2 = "Ứng lương tại LONG XUYÊN"
$ws.Range("B22").Value2 = 0
$ws.Range("A23").Value2 = "Tổng công tại SÓC TRĂNG"
$ws.Range("B23").Value2 = 0
$ws.Range("A24").Value2 = "Phụ cấp tại SÓC TRĂNG"
$ws.Range("B24").Value2 = 0
$ws.Range("A25").Value2 = "Lương công tác tại SÓC TRĂNG"
$ws.Range("B25").Value2 = 0
$ws.Range("A26").Value2 = "Lương cơ bản tại SÓC TRĂNG"
$ws.Range("B26").ClearContents()
$ws.Range("A27").Value2 = "Chiết khấu sale chính tại SÓC TRĂNG"
$ws.Range("B27").Value2 = 0
$ws.Range("A28").Value2 = "Chiết khấu sale phụ tại SÓC TRĂNG"
$ws.Range("B28").Value2 = 0
$ws.Range("A29").Value2 = "Đơn 1 bác sĩ tại SÓC TRĂNG"
$ws.Range("B29").Value2 = 0
$ws.Range("A30").Value2 = "Đơn 2 bác sĩ tại SÓC TRĂNG"
$ws.Range("B30").Value2 = 0
$ws.Range("A31").Value2 = "Công phụ phẫu 1 tại SÓC TRĂNG"
$ws.Range("B31").Value2 = 0
$ws.Range("A32").Value2 = "Công phụ phẫu 2 tại SÓC TRĂNG"
$ws.Range("B32").Value2 = 0
$ws.Range("A33").Value2 = "Ứng lương tại SÓC TRĂNG"
$ws.Range("B33").Value2 = 0
$ws.Range("A34").Value2 = "Tổng lương tại CẦN THƠ"
$ws.Range("B34").Value2 = 6980000
$ws.Range("A35").Value2 = "Tổng lương tại LONG XUYÊN"
$ws.Range("B35").Value2 = 0
$ws.Range("A36").Value2 = "Tổng lương tại SÓC TRĂNG"
$ws.Range("B36").Value2 = 0
$ws.Range("A37").Value2 = "Tổng lương"
$ws.Range("B37").Value2 = 6980000
